$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns(1).ColumnWidth = 10.81640625
